# Updated_Format_1 (Text Wrapper+Coloana B mai mica)
#
# Applies:
#  1. A3 label gets a 2-space prefix ("  Numele și Prenumele:")
#  2. The "index number" font (column A counters) recolors FF808080 -> FFE3E3E3
#     and gets centered (horizontal + vertical) alignment.
#  3. New "sentence count" values are written into column A for every word
#     block (rows 7,19,31,43,55,67,79,91,103,115,127,139), using the same
#     style as the other index-number cells.
#  4. Word-header cells (col A, fontId3/fill5) get horizontal+vertical center
#     plus wrap text.
#  5. Every other populated column-B cell (definition / example rows) gets
#     wrap text turned on, keeping its existing alignment/fill/font.
#  6. Column B width shrinks from 250.7109375 to 200.7109375.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Label text gets a two-space indent
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = "  Numele și Prenumele:"

# ---------------------------------------------------------------------------
# 2 & 3. Column A "number" cells: existing index numbers + new sentence counts
# ---------------------------------------------------------------------------
$newCounts = @{
    7   = 7
    19  = 10
    31  = 10
    43  = 10
    55  = 9
    67  = 8
    79  = 10
    91  = 5
    103 = 5
    115 = 10
    127 = 10
    139 = 10
}

foreach ($r in $newCounts.Keys) {
    $ws.Cells.Item($r, 1).Value2 = $newCounts[$r]
}

$existingNumberRows = @(1, 6, 18, 30, 42, 54, 66, 78, 90, 102, 114, 126, 138)
$numberRows = $existingNumberRows + @($newCounts.Keys)

# Apply one property at a time across the whole group - this lets same-shaped
# cells converge on a single shared style record instead of fragmenting into
# one new record per cell.
foreach ($r in $numberRows) {
    $ws.Cells.Item($r, 1).HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
}
foreach ($r in $numberRows) {
    $ws.Cells.Item($r, 1).VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
}
foreach ($r in $numberRows) {
    $ws.Cells.Item($r, 1).Font.Color = 14935011   # RGB(227,227,227) = FFE3E3E3
}

# ---------------------------------------------------------------------------
# 4. Word-header cells in column A: center align + wrap text
# ---------------------------------------------------------------------------
$wordHeaderRows = @(5, 17, 29, 41, 53, 65, 77, 89, 101, 113, 125, 137)
foreach ($r in $wordHeaderRows) {
    $ws.Cells.Item($r, 1).HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
}
foreach ($r in $wordHeaderRows) {
    $ws.Cells.Item($r, 1).VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
}
foreach ($r in $wordHeaderRows) {
    $ws.Cells.Item($r, 1).WrapText = $true
}

# ---------------------------------------------------------------------------
# 5. Every other already-formatted column-B cell gets wrap text
# ---------------------------------------------------------------------------
$wrapOnlyRows = @(
    5, 17, 29, 41, 53, 65, 77, 89, 101, 113, 125, 137,
    6, 8, 54, 66, 90, 92, 94, 102, 104, 106,
    7, 67, 91, 93, 103, 105,
    9, 11, 13, 15, 19, 21, 23, 25, 27, 31, 33, 35, 37, 39, 43, 45, 47, 49, 51,
    55, 57, 59, 61, 63, 69, 71, 73, 75, 79, 81, 83, 85, 87, 95, 97, 99, 107,
    109, 111, 115, 117, 119, 121, 123, 127, 129, 131, 133, 135, 139, 141, 143,
    145, 147,
    10, 12, 14, 18, 20, 22, 24, 26, 30, 32, 34, 36, 38, 42, 44, 46, 48, 50, 56,
    58, 60, 62, 68, 70, 72, 74, 78, 80, 82, 84, 86, 96, 98, 108, 110, 114, 116,
    118, 120, 122, 126, 128, 130, 132, 134, 138, 140, 142, 144, 146
) | Select-Object -Unique

foreach ($r in $wrapOnlyRows) {
    $ws.Cells.Item($r, 2).WrapText = $true
}

# ---------------------------------------------------------------------------
# 6. Column B narrower
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 200.7109375

Write-Output "done"
